$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows to append (columns A,B,C,D,E,F,G,H), matching the layout of
# the existing rows (Principle, Start Principle, BuyPrice, SellPrice,
# IsShortSell, Price Change %, Date, Profitable).
$rows = @(
    @(9816.16,             9822.0499999999993, 283.47000000000003, 283.29000000000002, $false, -0.06,              42613.7655787037,   $false),
    @(9764.1299999999992,  9816.16,             282.39,             280.89,             $false, -0.53,              42614.672962962963, $false),
    @(9792.4500000000007,  9764.1299999999992,  280.62,             281.44,             $false, 0.28999999999999998, 42615.750185185185, $true)
)

$startRow = 11
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]

    # Column G uses a date display format elsewhere in the sheet (row above
    # has the proper style); copy that formatting onto the new date cell
    # instead of creating a brand-new style entry.
    $ws.Range("G" + ($r - 1)).Copy() | Out-Null
    $ws.Range("G" + $r).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = $false
